$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; this pushes the existing rows 71-77 down to 72-78
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly price record
$ws.Cells.Item(71, 1).Value = 7
$ws.Cells.Item(71, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(71, 3).Value = "Ñuble"
$ws.Cells.Item(71, 4).Value = 45131
$ws.Cells.Item(71, 5).Value = 16
$ws.Cells.Item(71, 6).Value = 100112026
$ws.Cells.Item(71, 7).Value = "Haba"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 40
$ws.Cells.Item(71, 11).Value = 16000
$ws.Cells.Item(71, 12).Value = 16000
$ws.Cells.Item(71, 13).Value = 16000
$ws.Cells.Item(71, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(71, 16).Value = 640
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
